# Auto-generated Excel COM-interop script
# Adds Corequisites / Concurrent / Recommended columns (D, E, F) and
# shifts 'Terms Typically Offered' to column G, splitting out embedded
# corequisite/recommended text that had been concatenated into other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 4).Value = 'Corequisites'
$ws.Cells.Item(1, 5).Value = 'Concurrent'
$ws.Cells.Item(1, 6).Value = 'Recommended'
$ws.Cells.Item(1, 7).Value = 'Terms Typically Offered'

# --- Data rows (2-95) ---
# Row 2
$ws.Cells.Item(2, 3).Value = 'NA'
$ws.Cells.Item(2, 4).Value = 'NA'
$ws.Cells.Item(2, 5).Value = 'NA'
$ws.Cells.Item(2, 6).Value = 'NA'
$ws.Cells.Item(2, 7).Value = 'F, W, SP'

# Row 3
$ws.Cells.Item(3, 3).Value = 'MU 101 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(3, 4).Value = 'NA'
$ws.Cells.Item(3, 5).Value = 'NA'
$ws.Cells.Item(3, 6).Value = 'NA'
$ws.Cells.Item(3, 7).Value = 'F, W'

# Row 4
$ws.Cells.Item(4, 3).Value = 'MU 101; Music majors may be concurrently enrolled in MU 101 and MU 104.'
$ws.Cells.Item(4, 4).Value = 'NA'
$ws.Cells.Item(4, 5).Value = 'NA'
$ws.Cells.Item(4, 6).Value = 'NA'
$ws.Cells.Item(4, 7).Value = 'F'

# Row 5
$ws.Cells.Item(5, 3).Value = 'MU 103 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(5, 4).Value = 'NA'
$ws.Cells.Item(5, 5).Value = 'NA'
$ws.Cells.Item(5, 6).Value = 'NA'
$ws.Cells.Item(5, 7).Value = 'W, SP'

# Row 6
$ws.Cells.Item(6, 3).Value = 'MU 104 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(6, 4).Value = 'MU 103.'
$ws.Cells.Item(6, 5).Value = 'NA'
$ws.Cells.Item(6, 6).Value = 'NA'
$ws.Cells.Item(6, 7).Value = 'W '

# Row 7
$ws.Cells.Item(7, 3).Value = 'MU 106 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(7, 4).Value = 'NA'
$ws.Cells.Item(7, 5).Value = 'NA'
$ws.Cells.Item(7, 6).Value = 'NA'
$ws.Cells.Item(7, 7).Value = 'SP'

# Row 8
$ws.Cells.Item(8, 3).Value = 'MU 101 or consent of instructor.'
$ws.Cells.Item(8, 4).Value = 'NA'
$ws.Cells.Item(8, 5).Value = 'NA'
$ws.Cells.Item(8, 6).Value = 'NA'
$ws.Cells.Item(8, 7).Value = 'TBD'

# Row 9
$ws.Cells.Item(9, 3).Value = 'NA'
$ws.Cells.Item(9, 4).Value = 'NA'
$ws.Cells.Item(9, 5).Value = 'NA'
$ws.Cells.Item(9, 6).Value = 'NA'
$ws.Cells.Item(9, 7).Value = 'F, W, SP'

# Row 10
$ws.Cells.Item(10, 3).Value = 'Music major, minor, or consent of instructor.'
$ws.Cells.Item(10, 4).Value = 'NA'
$ws.Cells.Item(10, 5).Value = 'NA'
$ws.Cells.Item(10, 6).Value = 'NA'
$ws.Cells.Item(10, 7).Value = 'SP'

# Row 11
$ws.Cells.Item(11, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(11, 4).Value = 'NA'
$ws.Cells.Item(11, 5).Value = 'NA'
$ws.Cells.Item(11, 6).Value = 'NA'
$ws.Cells.Item(11, 7).Value = 'TBD'

# Row 12
$ws.Cells.Item(12, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(12, 4).Value = 'NA'
$ws.Cells.Item(12, 5).Value = 'NA'
$ws.Cells.Item(12, 6).Value = 'NA'
$ws.Cells.Item(12, 7).Value = 'F, W, SP'

# Row 13
$ws.Cells.Item(13, 3).Value = 'NA'
$ws.Cells.Item(13, 4).Value = 'NA'
$ws.Cells.Item(13, 5).Value = 'NA'
$ws.Cells.Item(13, 6).Value = 'NA'
$ws.Cells.Item(13, 7).Value = 'F, W, SP'

# Row 14
$ws.Cells.Item(14, 3).Value = 'MU 151 or consent of instructor; for non-music majors.'
$ws.Cells.Item(14, 4).Value = 'NA'
$ws.Cells.Item(14, 5).Value = 'NA'
$ws.Cells.Item(14, 6).Value = 'NA'
$ws.Cells.Item(14, 7).Value = 'F, W, SP'

# Row 15
$ws.Cells.Item(15, 3).Value = 'NA'
$ws.Cells.Item(15, 4).Value = 'NA'
$ws.Cells.Item(15, 5).Value = 'NA'
$ws.Cells.Item(15, 6).Value = 'NA'
$ws.Cells.Item(15, 7).Value = 'F, W, SP'

# Row 16
$ws.Cells.Item(16, 3).Value = 'NA'
$ws.Cells.Item(16, 4).Value = 'NA'
$ws.Cells.Item(16, 5).Value = 'NA'
$ws.Cells.Item(16, 6).Value = 'NA'
$ws.Cells.Item(16, 7).Value = 'F, W, SP'

# Row 17
$ws.Cells.Item(17, 3).Value = 'Music major or consent of instructor.'
$ws.Cells.Item(17, 4).Value = 'NA'
$ws.Cells.Item(17, 5).Value = 'NA'
$ws.Cells.Item(17, 6).Value = 'NA'
$ws.Cells.Item(17, 7).Value = 'F'

# Row 18
$ws.Cells.Item(18, 3).Value = 'MU 161 or consent of instructor.'
$ws.Cells.Item(18, 4).Value = 'NA'
$ws.Cells.Item(18, 5).Value = 'NA'
$ws.Cells.Item(18, 6).Value = 'NA'
$ws.Cells.Item(18, 7).Value = 'W'

# Row 19
$ws.Cells.Item(19, 3).Value = 'MU 162 or consent of instructor.'
$ws.Cells.Item(19, 4).Value = 'NA'
$ws.Cells.Item(19, 5).Value = 'NA'
$ws.Cells.Item(19, 6).Value = 'NA'
$ws.Cells.Item(19, 7).Value = 'SP'

# Row 20
$ws.Cells.Item(20, 3).Value = 'NA'
$ws.Cells.Item(20, 4).Value = 'MU 253 or piano topic in any of the following MU 150, MU 250, MU 350 or MU 450; or consent of instructor.'
$ws.Cells.Item(20, 5).Value = 'NA'
$ws.Cells.Item(20, 6).Value = 'NA'
$ws.Cells.Item(20, 7).Value = 'F, W, SP'

# Row 21
$ws.Cells.Item(21, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(21, 4).Value = 'NA'
$ws.Cells.Item(21, 5).Value = 'NA'
$ws.Cells.Item(21, 6).Value = 'NA'
$ws.Cells.Item(21, 7).Value = 'F, W, SP'

# Row 22
$ws.Cells.Item(22, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(22, 4).Value = 'NA'
$ws.Cells.Item(22, 5).Value = 'NA'
$ws.Cells.Item(22, 6).Value = 'NA'
$ws.Cells.Item(22, 7).Value = 'F, W, SP'

# Row 23
$ws.Cells.Item(23, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(23, 4).Value = 'NA'
$ws.Cells.Item(23, 5).Value = 'NA'
$ws.Cells.Item(23, 6).Value = 'NA'
$ws.Cells.Item(23, 7).Value = 'F, W, SP'

# Row 24
$ws.Cells.Item(24, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(24, 4).Value = 'NA'
$ws.Cells.Item(24, 5).Value = 'NA'
$ws.Cells.Item(24, 6).Value = 'NA'
$ws.Cells.Item(24, 7).Value = 'F, W, SP'

# Row 25
$ws.Cells.Item(25, 3).Value = 'Consent of instructor, based on audition.'
$ws.Cells.Item(25, 4).Value = 'NA'
$ws.Cells.Item(25, 5).Value = 'NA'
$ws.Cells.Item(25, 6).Value = 'NA'
$ws.Cells.Item(25, 7).Value = 'F, W, SP'

# Row 26
$ws.Cells.Item(26, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(26, 4).Value = 'NA'
$ws.Cells.Item(26, 5).Value = 'NA'
$ws.Cells.Item(26, 6).Value = 'NA'
$ws.Cells.Item(26, 7).Value = 'F, W'

# Row 27
$ws.Cells.Item(27, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(27, 4).Value = 'NA'
$ws.Cells.Item(27, 5).Value = 'NA'
$ws.Cells.Item(27, 6).Value = 'NA'
$ws.Cells.Item(27, 7).Value = 'F'

# Row 28
$ws.Cells.Item(28, 3).Value = 'NA'
$ws.Cells.Item(28, 4).Value = 'NA'
$ws.Cells.Item(28, 5).Value = 'NA'
$ws.Cells.Item(28, 6).Value = 'NA'
$ws.Cells.Item(28, 7).Value = 'F'

# Row 29
$ws.Cells.Item(29, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(29, 4).Value = 'NA'
$ws.Cells.Item(29, 5).Value = 'NA'
$ws.Cells.Item(29, 6).Value = 'NA'
$ws.Cells.Item(29, 7).Value = 'F, W, SP'

# Row 30
$ws.Cells.Item(30, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(30, 4).Value = 'NA'
$ws.Cells.Item(30, 5).Value = 'NA'
$ws.Cells.Item(30, 6).Value = 'NA'
$ws.Cells.Item(30, 7).Value = 'F, W, SP'

# Row 31
$ws.Cells.Item(31, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(31, 4).Value = 'NA'
$ws.Cells.Item(31, 5).Value = 'NA'
$ws.Cells.Item(31, 6).Value = 'NA'
$ws.Cells.Item(31, 7).Value = 'F, W, SP'

# Row 32
$ws.Cells.Item(32, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(32, 4).Value = 'NA'
$ws.Cells.Item(32, 5).Value = 'NA'
$ws.Cells.Item(32, 6).Value = 'NA'
$ws.Cells.Item(32, 7).Value = 'F, W, SP'

# Row 33
$ws.Cells.Item(33, 3).Value = 'NA'
$ws.Cells.Item(33, 4).Value = 'NA'
$ws.Cells.Item(33, 5).Value = 'NA'
$ws.Cells.Item(33, 6).Value = 'NA'
$ws.Cells.Item(33, 7).Value = 'F, W, SP'

# Row 34
$ws.Cells.Item(34, 3).Value = 'MU 150, MU 250, MU 350 or MU 450, or consent of instructor.'
$ws.Cells.Item(34, 4).Value = 'NA'
$ws.Cells.Item(34, 5).Value = 'NA'
$ws.Cells.Item(34, 6).Value = 'NA'
$ws.Cells.Item(34, 7).Value = 'F, W, SP'

# Row 35
$ws.Cells.Item(35, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(35, 4).Value = 'NA'
$ws.Cells.Item(35, 5).Value = 'NA'
$ws.Cells.Item(35, 6).Value = 'NA'
$ws.Cells.Item(35, 7).Value = 'F, W, SP'

# Row 36
$ws.Cells.Item(36, 3).Value = 'MU 108 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(36, 4).Value = 'MU 105.'
$ws.Cells.Item(36, 5).Value = 'NA'
$ws.Cells.Item(36, 6).Value = 'NA'
$ws.Cells.Item(36, 7).Value = 'F '

# Row 37
$ws.Cells.Item(37, 3).Value = 'MU 210 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(37, 4).Value = 'NA'
$ws.Cells.Item(37, 5).Value = 'NA'
$ws.Cells.Item(37, 6).Value = 'NA'
$ws.Cells.Item(37, 7).Value = 'W'

# Row 38
$ws.Cells.Item(38, 3).Value = 'MU 211 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(38, 4).Value = 'NA'
$ws.Cells.Item(38, 5).Value = 'NA'
$ws.Cells.Item(38, 6).Value = 'NA'
$ws.Cells.Item(38, 7).Value = 'SP'

# Row 39
$ws.Cells.Item(39, 3).Value = 'NA'
$ws.Cells.Item(39, 4).Value = 'NA'
$ws.Cells.Item(39, 5).Value = 'NA'
$ws.Cells.Item(39, 6).Value = 'NA'
$ws.Cells.Item(39, 7).Value = 'W, SP'

# Row 40
$ws.Cells.Item(40, 3).Value = 'NA'
$ws.Cells.Item(40, 4).Value = 'NA'
$ws.Cells.Item(40, 5).Value = 'NA'
$ws.Cells.Item(40, 6).Value = 'GE Area A1, GE Area A2, or GE Area A3.'
$ws.Cells.Item(40, 7).Value = 'F, W'

# Row 41
$ws.Cells.Item(41, 3).Value = 'NA'
$ws.Cells.Item(41, 4).Value = 'NA'
$ws.Cells.Item(41, 5).Value = 'NA'
$ws.Cells.Item(41, 6).Value = 'NA'
$ws.Cells.Item(41, 7).Value = 'TBD'

# Row 42
$ws.Cells.Item(42, 3).Value = '3 units of MU 150 and consent of instructor.'
$ws.Cells.Item(42, 4).Value = 'NA'
$ws.Cells.Item(42, 5).Value = 'NA'
$ws.Cells.Item(42, 6).Value = 'NA'
$ws.Cells.Item(42, 7).Value = 'F, W, SP'

# Row 43
$ws.Cells.Item(43, 3).Value = 'MU 153 or consent of instructor; for non-music majors.'
$ws.Cells.Item(43, 4).Value = 'NA'
$ws.Cells.Item(43, 5).Value = 'NA'
$ws.Cells.Item(43, 6).Value = 'NA'
$ws.Cells.Item(43, 7).Value = 'F, W, SP'

# Row 44
$ws.Cells.Item(44, 3).Value = 'Facility on a musical instrument or singing ability; MU 101 or consent of instructor.'
$ws.Cells.Item(44, 4).Value = 'NA'
$ws.Cells.Item(44, 5).Value = 'NA'
$ws.Cells.Item(44, 6).Value = 'NA'
$ws.Cells.Item(44, 7).Value = 'TBD'

# Row 45
$ws.Cells.Item(45, 3).Value = 'MU 163 or consent of instructor.'
$ws.Cells.Item(45, 4).Value = 'NA'
$ws.Cells.Item(45, 5).Value = 'NA'
$ws.Cells.Item(45, 6).Value = 'NA'
$ws.Cells.Item(45, 7).Value = 'F'

# Row 46
$ws.Cells.Item(46, 3).Value = 'MU 261 or consent of instructor.'
$ws.Cells.Item(46, 4).Value = 'NA'
$ws.Cells.Item(46, 5).Value = 'NA'
$ws.Cells.Item(46, 6).Value = 'NA'
$ws.Cells.Item(46, 7).Value = 'W'

# Row 47
$ws.Cells.Item(47, 3).Value = 'MU 262 or consent of instructor.'
$ws.Cells.Item(47, 4).Value = 'NA'
$ws.Cells.Item(47, 5).Value = 'NA'
$ws.Cells.Item(47, 6).Value = 'NA'
$ws.Cells.Item(47, 7).Value = 'SP'

# Row 48
$ws.Cells.Item(48, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(48, 4).Value = 'NA'
$ws.Cells.Item(48, 5).Value = 'NA'
$ws.Cells.Item(48, 6).Value = 'NA'
$ws.Cells.Item(48, 7).Value = 'W'

# Row 49
$ws.Cells.Item(49, 3).Value = 'Open to undergraduate students and consent of instructor.'
$ws.Cells.Item(49, 4).Value = 'NA'
$ws.Cells.Item(49, 5).Value = 'NA'
$ws.Cells.Item(49, 6).Value = 'NA'
$ws.Cells.Item(49, 7).Value = 'TBD'

# Row 50
$ws.Cells.Item(50, 3).Value = 'MU 105.'
$ws.Cells.Item(50, 4).Value = 'NA'
$ws.Cells.Item(50, 5).Value = 'NA'
$ws.Cells.Item(50, 6).Value = 'NA'
$ws.Cells.Item(50, 7).Value = 'TBD'

# Row 51
$ws.Cells.Item(51, 3).Value = 'MU 105 with a grade of C- or better, or consent of instructor.'
$ws.Cells.Item(51, 4).Value = 'NA'
$ws.Cells.Item(51, 5).Value = 'NA'
$ws.Cells.Item(51, 6).Value = 'NA'
$ws.Cells.Item(51, 7).Value = 'F'

# Row 52
$ws.Cells.Item(52, 3).Value = 'MU 303 or permission of instructor.'
$ws.Cells.Item(52, 4).Value = 'NA'
$ws.Cells.Item(52, 5).Value = 'NA'
$ws.Cells.Item(52, 6).Value = 'NA'
$ws.Cells.Item(52, 7).Value = 'W'

# Row 53
$ws.Cells.Item(53, 3).Value = 'MU 101, MU 120 or consent of instructor.'
$ws.Cells.Item(53, 4).Value = 'NA'
$ws.Cells.Item(53, 5).Value = 'NA'
$ws.Cells.Item(53, 6).Value = 'NA'
$ws.Cells.Item(53, 7).Value = 'F'

# Row 54
$ws.Cells.Item(54, 3).Value = 'MU 311 or permission of instructor.'
$ws.Cells.Item(54, 4).Value = 'NA'
$ws.Cells.Item(54, 5).Value = 'NA'
$ws.Cells.Item(54, 6).Value = 'NA'
$ws.Cells.Item(54, 7).Value = 'W'

# Row 55
$ws.Cells.Item(55, 3).Value = 'Completion of GE Area A1 with a grade of C- or better; and MU 105.'
$ws.Cells.Item(55, 4).Value = 'NA'
$ws.Cells.Item(55, 5).Value = 'NA'
$ws.Cells.Item(55, 6).Value = 'MU 120.'
$ws.Cells.Item(55, 7).Value = 'F '

# Row 56
$ws.Cells.Item(56, 3).Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one course in GE Area C.'
$ws.Cells.Item(56, 4).Value = 'NA'
$ws.Cells.Item(56, 5).Value = 'NA'
$ws.Cells.Item(56, 6).Value = 'NA'
$ws.Cells.Item(56, 7).Value = 'F, W, SP'

# Row 57
$ws.Cells.Item(57, 3).Value = 'MU 105.'
$ws.Cells.Item(57, 4).Value = 'NA'
$ws.Cells.Item(57, 5).Value = 'NA'
$ws.Cells.Item(57, 6).Value = 'MU 120.'
$ws.Cells.Item(57, 7).Value = 'SP '

# Row 58
$ws.Cells.Item(58, 3).Value = 'MU 121 or consent of instructor.'
$ws.Cells.Item(58, 4).Value = 'NA'
$ws.Cells.Item(58, 5).Value = 'NA'
$ws.Cells.Item(58, 6).Value = 'NA'
$ws.Cells.Item(58, 7).Value = 'SP'

# Row 59
$ws.Cells.Item(59, 3).Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one course in GE Area C.'
$ws.Cells.Item(59, 4).Value = 'NA'
$ws.Cells.Item(59, 5).Value = 'NA'
$ws.Cells.Item(59, 6).Value = 'NA'
$ws.Cells.Item(59, 7).Value = 'TBD'

# Row 60
$ws.Cells.Item(60, 3).Value = 'MU 320 or consent of instructor.'
$ws.Cells.Item(60, 4).Value = 'NA'
$ws.Cells.Item(60, 5).Value = 'NA'
$ws.Cells.Item(60, 6).Value = 'MU 120.'
$ws.Cells.Item(60, 7).Value = 'W '

# Row 61
$ws.Cells.Item(61, 3).Value = 'MU 320 or consent of instructor.'
$ws.Cells.Item(61, 4).Value = 'NA'
$ws.Cells.Item(61, 5).Value = 'NA'
$ws.Cells.Item(61, 6).Value = 'MU 120.'
$ws.Cells.Item(61, 7).Value = 'SP '

# Row 62
$ws.Cells.Item(62, 3).Value = 'MU 105.'
$ws.Cells.Item(62, 4).Value = 'NA'
$ws.Cells.Item(62, 5).Value = 'NA'
$ws.Cells.Item(62, 6).Value = 'NA'
$ws.Cells.Item(62, 7).Value = 'TBD'

# Row 63
$ws.Cells.Item(63, 3).Value = 'MU 105 or consent of instructor.'
$ws.Cells.Item(63, 4).Value = 'NA'
$ws.Cells.Item(63, 5).Value = 'NA'
$ws.Cells.Item(63, 6).Value = 'NA'
$ws.Cells.Item(63, 7).Value = 'F'

# Row 64
$ws.Cells.Item(64, 3).Value = 'MU 340.'
$ws.Cells.Item(64, 4).Value = 'NA'
$ws.Cells.Item(64, 5).Value = 'NA'
$ws.Cells.Item(64, 6).Value = 'NA'
$ws.Cells.Item(64, 7).Value = 'W'

# Row 65
$ws.Cells.Item(65, 3).Value = 'MU 340.'
$ws.Cells.Item(65, 4).Value = 'NA'
$ws.Cells.Item(65, 5).Value = 'NA'
$ws.Cells.Item(65, 6).Value = 'NA'
$ws.Cells.Item(65, 7).Value = 'SP'

# Row 66
$ws.Cells.Item(66, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(66, 4).Value = 'NA'
$ws.Cells.Item(66, 5).Value = 'NA'
$ws.Cells.Item(66, 6).Value = 'NA'
$ws.Cells.Item(66, 7).Value = 'F, W, SP'

# Row 67
$ws.Cells.Item(67, 3).Value = 'MU 105.'
$ws.Cells.Item(67, 4).Value = 'NA'
$ws.Cells.Item(67, 5).Value = 'NA'
$ws.Cells.Item(67, 6).Value = 'NA'
$ws.Cells.Item(67, 7).Value = 'TBD'

# Row 68
$ws.Cells.Item(68, 3).Value = 'MU 105.'
$ws.Cells.Item(68, 4).Value = 'NA'
$ws.Cells.Item(68, 5).Value = 'NA'
$ws.Cells.Item(68, 6).Value = 'NA'
$ws.Cells.Item(68, 7).Value = 'TBD'

# Row 69
$ws.Cells.Item(69, 3).Value = 'MU 101.'
$ws.Cells.Item(69, 4).Value = 'NA'
$ws.Cells.Item(69, 5).Value = 'NA'
$ws.Cells.Item(69, 6).Value = 'NA'
$ws.Cells.Item(69, 7).Value = 'TBD'

# Row 70
$ws.Cells.Item(70, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(70, 4).Value = 'NA'
$ws.Cells.Item(70, 5).Value = 'NA'
$ws.Cells.Item(70, 6).Value = 'NA'
$ws.Cells.Item(70, 7).Value = 'TBD'

# Row 71
$ws.Cells.Item(71, 3).Value = 'Junior standing and MU 168 or consent of instructor.'
$ws.Cells.Item(71, 4).Value = 'NA'
$ws.Cells.Item(71, 5).Value = 'NA'
$ws.Cells.Item(71, 6).Value = 'NA'
$ws.Cells.Item(71, 7).Value = 'F, W, SP'

# Row 72
$ws.Cells.Item(72, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(72, 4).Value = 'NA'
$ws.Cells.Item(72, 5).Value = 'NA'
$ws.Cells.Item(72, 6).Value = 'NA'
$ws.Cells.Item(72, 7).Value = 'F, W, SP'

# Row 73
$ws.Cells.Item(73, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(73, 4).Value = 'NA'
$ws.Cells.Item(73, 5).Value = 'NA'
$ws.Cells.Item(73, 6).Value = 'NA'
$ws.Cells.Item(73, 7).Value = 'F, W, SP'

# Row 74
$ws.Cells.Item(74, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(74, 4).Value = 'NA'
$ws.Cells.Item(74, 5).Value = 'NA'
$ws.Cells.Item(74, 6).Value = 'NA'
$ws.Cells.Item(74, 7).Value = 'F, W, SP'

# Row 75
$ws.Cells.Item(75, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(75, 4).Value = 'NA'
$ws.Cells.Item(75, 5).Value = 'NA'
$ws.Cells.Item(75, 6).Value = 'NA'
$ws.Cells.Item(75, 7).Value = 'F, W, SP'

# Row 76
$ws.Cells.Item(76, 3).Value = 'Junior standing and consent of instructor, based on audition.'
$ws.Cells.Item(76, 4).Value = 'NA'
$ws.Cells.Item(76, 5).Value = 'NA'
$ws.Cells.Item(76, 6).Value = 'NA'
$ws.Cells.Item(76, 7).Value = 'F, W, SP'

# Row 77
$ws.Cells.Item(77, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(77, 4).Value = 'NA'
$ws.Cells.Item(77, 5).Value = 'NA'
$ws.Cells.Item(77, 6).Value = 'NA'
$ws.Cells.Item(77, 7).Value = 'F, W'

# Row 78
$ws.Cells.Item(78, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(78, 4).Value = 'NA'
$ws.Cells.Item(78, 5).Value = 'NA'
$ws.Cells.Item(78, 6).Value = 'NA'
$ws.Cells.Item(78, 7).Value = 'F'

# Row 79
$ws.Cells.Item(79, 3).Value = 'MU 178 or consent of instructor.'
$ws.Cells.Item(79, 4).Value = 'NA'
$ws.Cells.Item(79, 5).Value = 'NA'
$ws.Cells.Item(79, 6).Value = 'NA'
$ws.Cells.Item(79, 7).Value = 'F'

# Row 80
$ws.Cells.Item(80, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(80, 4).Value = 'NA'
$ws.Cells.Item(80, 5).Value = 'NA'
$ws.Cells.Item(80, 6).Value = 'NA'
$ws.Cells.Item(80, 7).Value = 'F, W, SP'

# Row 81
$ws.Cells.Item(81, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(81, 4).Value = 'NA'
$ws.Cells.Item(81, 5).Value = 'NA'
$ws.Cells.Item(81, 6).Value = 'NA'
$ws.Cells.Item(81, 7).Value = 'F, W, SP'

# Row 82
$ws.Cells.Item(82, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(82, 4).Value = 'NA'
$ws.Cells.Item(82, 5).Value = 'NA'
$ws.Cells.Item(82, 6).Value = 'NA'
$ws.Cells.Item(82, 7).Value = 'F, W, SP'

# Row 83
$ws.Cells.Item(83, 3).Value = 'Junior standing and consent of instructor.'
$ws.Cells.Item(83, 4).Value = 'NA'
$ws.Cells.Item(83, 5).Value = 'NA'
$ws.Cells.Item(83, 6).Value = 'NA'
$ws.Cells.Item(83, 7).Value = 'F, W, SP'

# Row 84
$ws.Cells.Item(84, 3).Value = 'MU 188 or consent of instructor.'
$ws.Cells.Item(84, 4).Value = 'NA'
$ws.Cells.Item(84, 5).Value = 'NA'
$ws.Cells.Item(84, 6).Value = 'NA'
$ws.Cells.Item(84, 7).Value = 'F, W, SP'

# Row 85
$ws.Cells.Item(85, 3).Value = 'MU 150, MU 250, MU 350 or MU 450, or consent of instructor.'
$ws.Cells.Item(85, 4).Value = 'NA'
$ws.Cells.Item(85, 5).Value = 'NA'
$ws.Cells.Item(85, 6).Value = 'NA'
$ws.Cells.Item(85, 7).Value = 'F, W, SP'

# Row 86
$ws.Cells.Item(86, 3).Value = 'Junior standing and consent of department head.'
$ws.Cells.Item(86, 4).Value = 'NA'
$ws.Cells.Item(86, 5).Value = 'NA'
$ws.Cells.Item(86, 6).Value = 'NA'
$ws.Cells.Item(86, 7).Value = 'F, W, SP'

# Row 87
$ws.Cells.Item(87, 3).Value = 'MU 312.'
$ws.Cells.Item(87, 4).Value = 'NA'
$ws.Cells.Item(87, 5).Value = 'NA'
$ws.Cells.Item(87, 6).Value = 'NA'
$ws.Cells.Item(87, 7).Value = 'SP'

# Row 88
$ws.Cells.Item(88, 3).Value = 'MU 312.'
$ws.Cells.Item(88, 4).Value = 'NA'
$ws.Cells.Item(88, 5).Value = 'NA'
$ws.Cells.Item(88, 6).Value = 'NA'
$ws.Cells.Item(88, 7).Value = 'TBD'

# Row 89
$ws.Cells.Item(89, 3).Value = 'MU 303 and MU 320;'
$ws.Cells.Item(89, 4).Value = 'NA'
$ws.Cells.Item(89, 5).Value = 'NA'
$ws.Cells.Item(89, 6).Value = 'MU 120; or consent of instructor.'
$ws.Cells.Item(89, 7).Value = 'F '

# Row 90
$ws.Cells.Item(90, 3).Value = 'MU 303; and one of the MU 331, MU 332, or MU 431.'
$ws.Cells.Item(90, 4).Value = 'NA'
$ws.Cells.Item(90, 5).Value = 'NA'
$ws.Cells.Item(90, 6).Value = 'MU 305.'
$ws.Cells.Item(90, 7).Value = 'W '

# Row 91
$ws.Cells.Item(91, 3).Value = 'Consent of instructor.'
$ws.Cells.Item(91, 4).Value = 'NA'
$ws.Cells.Item(91, 5).Value = 'NA'
$ws.Cells.Item(91, 6).Value = 'NA'
$ws.Cells.Item(91, 7).Value = 'F, W, SP'

# Row 92
$ws.Cells.Item(92, 3).Value = 'Senior standing and consent of department head.'
$ws.Cells.Item(92, 4).Value = 'NA'
$ws.Cells.Item(92, 5).Value = 'NA'
$ws.Cells.Item(92, 6).Value = 'NA'
$ws.Cells.Item(92, 7).Value = 'F, W, SP'

# Row 93
$ws.Cells.Item(93, 3).Value = 'MU 341, or consent of instructor.'
$ws.Cells.Item(93, 4).Value = 'NA'
$ws.Cells.Item(93, 5).Value = 'NA'
$ws.Cells.Item(93, 6).Value = 'NA'
$ws.Cells.Item(93, 7).Value = 'TBD'

# Row 94
$ws.Cells.Item(94, 3).Value = 'MU 101 or MU 103; and MU 172 or MU 173 or MU 174.'
$ws.Cells.Item(94, 4).Value = 'NA'
$ws.Cells.Item(94, 5).Value = 'NA'
$ws.Cells.Item(94, 6).Value = 'NA'
$ws.Cells.Item(94, 7).Value = 'TBD'

# Row 95
$ws.Cells.Item(95, 3).Value = 'MU 331, MU 332, MU 431, MU 432, or consent of instructor.'
$ws.Cells.Item(95, 4).Value = 'NA'
$ws.Cells.Item(95, 5).Value = 'NA'
$ws.Cells.Item(95, 6).Value = 'NA'
$ws.Cells.Item(95, 7).Value = 'TBD'
